# Apply "repull data, push all data, mean calculation" updates to column F (dSF)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -2
    6  = 12
    8  = -2
    9  = 3
    18 = 1
    20 = -1
    21 = 0
    27 = -1
    28 = 1
    30 = -2
    39 = -4
    42 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
